{"js": "// The four \"<id>...</id>\" markers in this document are each split across\n// three runs: \"<id>\" (Courier New / #7f6000), the bare id text (plain\n// black run), and \"</id>\" (Courier New / #7f6000 again). The edit merges\n// each triplet into a single run carrying the \"<id>...</id>\" text using\n// the Courier New / #7f6000 formatting of the surrounding runs.\n\nconst body = context.document.body;\n\n// Find every literal \"<id>\" marker in the document body.\nconst idMarkers = body.search(\"<id>\", { matchCase: true });\nidMarkers.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < idMarkers.items.length; i++) {\n  const marker = idMarkers.items[i];\n\n  // The \"<id>\", id text and \"</id>\" runs all live in the same paragraph\n  // (nothing else shares that paragraph), so the paragraph's own range\n  // is exactly the \"<id>...</id>\" span we need to collapse into one run.\n  const para = marker.paragraphs.getFirst();\n  const fullRange = para.getRange();\n  fullRange.load(\"text\");\n  await context.sync();\n\n  const combinedText = fullRange.text; // e.g. \"<id>p076v_1</id>\"\n\n  // Re-inserting the same text as a single replace collapses the three\n  // runs into one run, inheriting the formatting of the first run\n  // (\"<id>\"'s Courier New / #7f6000), exactly like the target edit.\n  fullRange.insertText(combinedText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The \"<id>...</id>\" markers in this document are each split across three\n# runs: \"<id>\" (Courier New / #7f6000), the bare id text (plain black\n# run), and \"</id>\" (Courier New / #7f6000 again) \u2014 all three runs alone\n# in their own paragraph. The edit merges each triplet into a single run\n# carrying the full \"<id>...</id>\" text, using the Courier New / #7f6000\n# formatting of the surrounding \"<id>\"/\"</id>\" runs.\n\n$d = $word.ActiveDocument\n\n# Discover every \"<id>...</id>\" paragraph and capture its full text\n# (rather than hardcoding the id values), so the fix-up below applies to\n# whatever ids are actually present.\n$targets = New-Object System.Collections.ArrayList\nforeach ($para in $d.Paragraphs) {\n  $t = $para.Range.Text\n  $trimmed = $t.TrimEnd([char]13, [char]7)\n  if ($trimmed.StartsWith(\"<id>\") -and $trimmed.EndsWith(\"</id>\")) {\n    [void]$targets.Add($trimmed)\n  }\n}\n\nforeach ($full in $targets) {\n  # The text between \"<id>\" and \"</id>\" uniquely locates this paragraph\n  # (searching for the generic \"<id>\" marker would always re-hit the\n  # first, already-fixed occurrence).\n  $inner = $full.Substring(4, $full.Length - 9)\n\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.Text = $inner\n  $find.Execute() | Out-Null\n\n  if ($find.Found) {\n    # The \"<id>\", id text and \"</id>\" runs all live in one paragraph by\n    # themselves, so the paragraph range (minus its trailing paragraph\n    # mark) is exactly the \"<id>...</id>\" span to collapse into one run.\n    $para = $rng.Paragraphs(1)\n    $prange = $para.Range\n    $prange.MoveEnd(1, -1) | Out-Null\n\n    # Assigning identical text back is a no-op, so first stamp a\n    # placeholder (it inherits the \"<id>\" run's Courier New / #7f6000\n    # formatting), then write the real combined text. This collapses the\n    # three original runs into a single run, matching the target edit.\n    $prange.Text = \"ZZPLACEHOLDERZZ\"\n\n    $prange2 = $para.Range\n    $prange2.MoveEnd(1, -1) | Out-Null\n    $prange2.Text = $full\n  }\n}\n"}
